# The crawl was re-run later the same day (2022-08-03 20:59:25 instead of
# 07:00:37), and two products that dropped out of the scraped result set
# need to be removed from the sheet:
#   - row 3: id 4947421 "Oecoplan Taschentuch Calendula Box 3x  80ST"
#   - row 5: id 6695141 "Prix Garantie feuchtes Toilettenpapier 2x70 Stück"
# Deleting them shifts every following row up by one (or two), and every
# remaining row's timestamp (column O) is refreshed to the new crawl time.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the higher row index first so row 3's index doesn't shift before
# we get to it.
$ws.Rows.Item(5).Delete()
$ws.Rows.Item(3).Delete()

# Refresh the timestamp column (O = column 15) for the meta row (row 2)
# and every remaining product row.
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 15).Value = "2022-08-03 20:59:25"
}
